# CSwap Deliverable 2 - Sprint Backlog update
# Updates "Story Points Completed" (column H) on the "Sprint 1" sheet,
# moves the active selection, and nudges a handful of column widths that
# were resized in the same save.

$wb = $excel.ActiveWorkbook

$sprint1 = $wb.Worksheets.Item("Sprint 1")
$sprint2 = $wb.Worksheets.Item("Sprint 2")

# --- Story Points Completed updates (column H) -----------------------
$sprint1.Range("H2").Value = 2
$sprint1.Range("H3").Value = 2
$sprint1.Range("H4").Value = 0.5
$sprint1.Range("H5").Value = 0.5
$sprint1.Range("H6").Value = 0.5

# --- Column width tweaks on "Sprint 1" --------------------------------
# (ColumnWidth is in "characters"; the stored OOXML width includes a
# fixed ~5/6 character padding offset, so subtract that to land on the
# target stored width.)
$offset = 5 / 6
$sprint1.Columns.Item(2).ColumnWidth = 81.58203125 - $offset   # B
$sprint1.Columns.Item(3).ColumnWidth = 35.83203125 - $offset   # C
$sprint1.Columns.Item(4).ColumnWidth = 13 - $offset            # D
$sprint1.Columns.Item(6).ColumnWidth = 9.58203125 - $offset    # F
$sprint1.Columns.Item(8).ColumnWidth = 19.6640625 - $offset    # H
$sprint1.Columns.Item(9).ColumnWidth = 30.5 - $offset          # I

# --- Column width tweak on "Sprint 2" ---------------------------------
for ($i = 1; $i -le 26; $i++) {
    $sprint2.Columns.Item($i).ColumnWidth = 10.58203125 - $offset
}

# --- Selection moves to B8 on "Sprint 1" ------------------------------
$sprint1.Activate() | Out-Null
$sprint1.Range("B8").Select() | Out-Null
